$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.876.72"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "'2.437.78"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'560.04"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'162.06"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.514"
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("E9").Value = "  +11.06%  "
$ws.Range("D10").Value = "'0.163"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "'4.61"
$ws.Range("E12").Value = "  -5.03%  "
$ws.Range("E13").Value = "  +4.78%  "
$ws.Range("D14").Value = "'68.773.51"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "'2.887.29"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "'2.440.27"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D19").Value = "'338.66"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "'6.94"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'67.07"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("D26").Value = "'2.565.45"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").Value = "'8.19"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").Value = "'0.0₃0817"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "'7.12"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "'427.68"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("D35").Value = "'160.45"
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("D36").Value = "'19.02"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").Value = "'1.50"
$ws.Range("E41").Value = "  +3.34%  "
$ws.Range("D42").Value = "'4.34"
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "'2.03"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("D45").Value = "'3.34"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "'130.10"
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").Value = "'0.0719"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("D48").Value = "'0.481"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("D50").Value = "'0.0923"
$ws.Range("E50").Value = "  +1.48%  "
